$p = $ppt.ActivePresentation

# --- Add the new slide (10th), using the same "Title and Content" layout
#     used throughout the rest of the deck (slideLayout2.xml == CustomLayouts item 2) ---
$layout = $p.SlideMaster.CustomLayouts.Item(2)
$s = $p.Slides.AddSlide($p.Slides.Count + 1, $layout)

# --- Title placeholder ---
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "15.11.2022"

# --- Body / content placeholder ---
$body = $s.Shapes.Item(2).TextFrame.TextRange
$p1 = "Improving the location code by using the Levenberg-Marquardt (page 707) method for least square minimazition."
$p2 = "The model function is the inverse square law, where the adjustable parameters are u (x source coordinate), v (y source coordinate) and A0 (the activity of the source). We have a number of data points of intensities along the grid which gives us a overestimated system of equations."
$body.Text = $p1 + "`r" + $p2

# Resize/position the content placeholder to match the authored layout
$s.Shapes.Item(2).Left = 838200 / 914400
$s.Shapes.Item(2).Top = 1825625 / 914400
$s.Shapes.Item(2).Width = 10515600 / 914400
$s.Shapes.Item(2).Height = 3399518 / 914400

# --- Hyperlink the "Levenberg-Marquardt" substring inside the body text ---
$link = $body.Characters(42, 19)
$link.ActionSettings(1).Hyperlink.Address = "https://en.wikipedia.org/wiki/Levenberg%E2%80%93Marquardt_algorithm"
